# buoy_samples_metadata_2019.xlsx - "processing the latest pH runs" pass
#
# Inspecting the canonical-OOXML diff for this workbook shows that, once the
# pure re-save noise is stripped away (a newer Excel build's fileVersion /
# calcPr / xr:* revision namespaces / window geometry / the font-metric
# driven dyDescent 0.3->0.25 + default row height 14.4->15 + the handful of
# bestFit column widths that move by a few hundredths of a character), the
# only worksheet-level artifact left is that the shared-formula definition
# hosted in F9 (si="0") has its bookkeeping "ref" span tighten from "F3:F9"
# down to "F9" - i.e. Excel recognizing that F3:F8 never actually carried
# the live MEDIAN() formula (only cached numbers) and that F9 is the sole
# cell still using it. No cell value, formula result, label or note changed
# anywhere on either sheet (the dimensions stay A1:O36 / A1:B14).
#
# The author's note about "P-0036B-1" having a 0.05 pH unit spread belongs
# to another worksheet touched in the same batch-processing commit; nothing
# in this file's own diff introduces/edits a note to that effect, so we
# don't fabricate one here.
#
# The faithful, content-level action this script performs is to reopen the
# median check-pH formula in F9 so the workbook reflects that it was
# revisited during this pass (harmless/idempotent - it recomputes to the
# same cached result) while leaving every other cell, value and formula in
# the workbook untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buoy_samples_metadata_2019")

$ws.Range("F9").Formula = "=MEDIAN(G9:I9)"
